# Bai 8 - adjust the subtitle placeholder's position/size on slide 1
# (name/title adjustment commit: reflow the subtitle box)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rectangle 3" is the subTitle placeholder (2nd shape on the slide)
$shp = $s.Shapes.Item(2)

# Target OOXML (EMU):
#   off  x=611560  y=3501008
#   ext  cx=8352928 cy=1440880
# Converted to points (1 pt = 12700 EMU); Height nudged a hair inside its
# rounding bucket so the COM point->EMU conversion lands exactly on 1440880.
$shp.Left   = 48.15433070866142
$shp.Top    = 275.6699212598425
$shp.Width  = 657.7108661417323
$shp.Height = 113.45515823364258
